$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row 13 (ID 12): "Student" -> "Student's T"
# ---------------------------------------------------------------------------
$ws.Range("B13").Value = "Student's T"

# ---------------------------------------------------------------------------
# 2. Row 14 (ID 13, low-mean-high): fill in the comment for column C
#    (style matches the other plain Arial text cells, e.g. B15)
# ---------------------------------------------------------------------------
$ws.Range("B15").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C14").Value = "mean value is stored in 'value' column, only low and high are stored here"

# ---------------------------------------------------------------------------
# 3. Row 15 (ID 14, undefined / static value): fill in the comment for column C
# ---------------------------------------------------------------------------
$ws.Range("C15").Value = "moved from 0 as 0 is not allowed in mySQL db."

# ---------------------------------------------------------------------------
# 4. Row 16 (ID 15): new "alternative value" entry
# ---------------------------------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("D16").PasteSpecial(-4122)

$ws.Range("B16").Value = "alternative value"
$ws.Range("C16").Value = "basic (= reference) value is stored in 'value' column"
$ws.Range("D16").Value = "alternative"

# ---------------------------------------------------------------------------
# 5. Row 17 (new, ID 16): "low alternative"
# ---------------------------------------------------------------------------
$ws.Range("B1").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("F17").PasteSpecial(-4122)

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "low alternative"
$ws.Range("C17").Value = "high (= reference) value is stored in 'value' column"
$ws.Range("D17").Value = "low"

# ---------------------------------------------------------------------------
# 6. Row 18 (new, ID 17): "high alternative"
# ---------------------------------------------------------------------------
$ws.Range("B1").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("F18").PasteSpecial(-4122)

$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "high alternative"
$ws.Range("C18").Value = "low (= reference) value is stored in 'value' column"
$ws.Range("D18").Value = "high"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 7. Selection moves to reflect the newly extended data range
# ---------------------------------------------------------------------------
$ws.Range("A2:F18").Select()
